{"js": "// Replace the date line and every \"a\u00f7b=\" exercise in the table with the\n// new values from the commit. Each old value occurs exactly once in the\n// document, but a handful of *new* values collide with *old* values that\n// occur elsewhere (e.g. \"73\u00f73=\" -> \"92\u00f79=\", while \"92\u00f79=\" is itself an old\n// value elsewhere). To avoid double-touching text, we first locate every\n// target range (by searching for the distinct old strings against the\n// pristine document) and only then perform the text replacements.\nconst mapping = [\n  [\"2024-03-21 Thursday\", \"2024-03-22 Friday\"],\n  [\"75\u00f79=\", \"33\u00f72=\"],\n  [\"99\u00f79=\", \"79\u00f73=\"],\n  [\"95\u00f73=\", \"69\u00f74=\"],\n  [\"74\u00f77=\", \"87\u00f79=\"],\n  [\"83\u00f75=\", \"55\u00f72=\"],\n  [\"14\u00f76=\", \"53\u00f72=\"],\n  [\"92\u00f79=\", \"52\u00f73=\"],\n  [\"65\u00f76=\", \"68\u00f75=\"],\n  [\"53\u00f74=\", \"95\u00f79=\"],\n  [\"30\u00f72=\", \"89\u00f73=\"],\n  [\"71\u00f78=\", \"80\u00f72=\"],\n  [\"41\u00f77=\", \"61\u00f79=\"],\n  [\"21\u00f74=\", \"24\u00f73=\"],\n  [\"80\u00f75=\", \"36\u00f73=\"],\n  [\"73\u00f73=\", \"92\u00f79=\"],\n  [\"63\u00f74=\", \"21\u00f75=\"],\n  [\"96\u00f74=\", \"19\u00f73=\"],\n  [\"30\u00f77=\", \"68\u00f74=\"],\n  [\"97\u00f77=\", \"28\u00f76=\"],\n  [\"18\u00f77=\", \"61\u00f77=\"],\n  [\"32\u00f77=\", \"27\u00f77=\"],\n  [\"83\u00f73=\", \"22\u00f73=\"],\n  [\"22\u00f75=\", \"41\u00f72=\"],\n  [\"74\u00f76=\", \"77\u00f75=\"],\n  [\"55\u00f74=\", \"59\u00f79=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: search for each old value (exact, case-sensitive, no wildcards)\n// and capture the resulting RangeCollection before mutating anything.\nconst pending = mapping.map(([oldText, newText]) => {\n  const results = body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n    ignoreSpace: false,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  return { newText, results };\n});\n\nawait context.sync();\n\n// Phase 2: replace the text in-place for each found range. insertText with\n// the \"Replace\" location swaps just the matched text, keeping the run's\n// formatting (font, size, paragraph alignment) untouched.\nfor (const { newText, results } of pending) {\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found while applying edits: \" + newText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every \"a\u00f7b=\" exercise to the new values\n# from the commit. Find/Replace (wdReplaceAll) swaps only the matched text\n# run, leaving paragraph/run formatting (font, size, alignment) untouched.\n#\n# NOTE ON ORDERING: every old value is unique and occurs exactly once, but\n# one *new* value (\"92\u00f79=\", produced by \"73\u00f73=\") collides with an *old*\n# value that appears elsewhere in the table (\"92\u00f79=\" -> \"52\u00f73=\"). If that\n# old \"92\u00f79=\" were replaced after the new one was written, the Find would\n# re-match the freshly written text and clobber it. So the \"92\u00f79=\" ->\n# \"52\u00f73=\" replacement is ordered before the \"73\u00f73=\" -> \"92\u00f79=\" replacement\n# below (all other pairs are independent and keep the diff's natural\n# order).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Exact($doc, [string]$oldText, [string]$newText) {\n    $range = $doc.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    # wdFindContinue=1 for Wrap, wdReplaceAll=2 for Replace \u2014 only one\n    # occurrence of each search string exists at call time, so ReplaceAll\n    # vs ReplaceOne makes no difference here.\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nReplace-Exact $d \"2024-03-21 Thursday\" \"2024-03-22 Friday\"\n\nReplace-Exact $d \"75\u00f79=\" \"33\u00f72=\"\nReplace-Exact $d \"99\u00f79=\" \"79\u00f73=\"\nReplace-Exact $d \"95\u00f73=\" \"69\u00f74=\"\nReplace-Exact $d \"74\u00f77=\" \"87\u00f79=\"\nReplace-Exact $d \"83\u00f75=\" \"55\u00f72=\"\n\nReplace-Exact $d \"14\u00f76=\" \"53\u00f72=\"\nReplace-Exact $d \"92\u00f79=\" \"52\u00f73=\"\nReplace-Exact $d \"73\u00f73=\" \"92\u00f79=\"\nReplace-Exact $d \"65\u00f76=\" \"68\u00f75=\"\nReplace-Exact $d \"53\u00f74=\" \"95\u00f79=\"\nReplace-Exact $d \"30\u00f72=\" \"89\u00f73=\"\n\nReplace-Exact $d \"71\u00f78=\" \"80\u00f72=\"\nReplace-Exact $d \"41\u00f77=\" \"61\u00f79=\"\nReplace-Exact $d \"21\u00f74=\" \"24\u00f73=\"\nReplace-Exact $d \"80\u00f75=\" \"36\u00f73=\"\n\nReplace-Exact $d \"63\u00f74=\" \"21\u00f75=\"\nReplace-Exact $d \"96\u00f74=\" \"19\u00f73=\"\nReplace-Exact $d \"30\u00f77=\" \"68\u00f74=\"\nReplace-Exact $d \"97\u00f77=\" \"28\u00f76=\"\nReplace-Exact $d \"18\u00f77=\" \"61\u00f77=\"\n\nReplace-Exact $d \"32\u00f77=\" \"27\u00f77=\"\nReplace-Exact $d \"83\u00f73=\" \"22\u00f73=\"\nReplace-Exact $d \"22\u00f75=\" \"41\u00f72=\"\nReplace-Exact $d \"74\u00f76=\" \"77\u00f75=\"\nReplace-Exact $d \"55\u00f74=\" \"59\u00f79=\"\n"}
